$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $origStyle = $Range.Style
    $Range.Value = "'" + $Text
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '46.031.94'
Set-TextValue $ws.Range("E2") '  -1.68%  '

Set-TextValue $ws.Range("D3") '2.379.31'
Set-TextValue $ws.Range("E3") '  +3.29%  '

Set-TextValue $ws.Range("E4") '  -0.02%  '

Set-TextValue $ws.Range("D5") '300.93'
Set-TextValue $ws.Range("E5") '  -0.72%  '

Set-TextValue $ws.Range("D6") '98.84'
Set-TextValue $ws.Range("E6") '  -2.70%  '

Set-TextValue $ws.Range("D7") '0.563'
Set-TextValue $ws.Range("E7") '  -1.22%  '

Set-TextValue $ws.Range("E8") '  +0.06%  '

Set-TextValue $ws.Range("E9") '  -4.47%  '

Set-TextValue $ws.Range("D10") '34.56'
Set-TextValue $ws.Range("E10") '  -6.23%  '

Set-TextValue $ws.Range("D11") '0.0789'
Set-TextValue $ws.Range("E11") '  -2.08%  '

Set-TextValue $ws.Range("E12") '  -4.52%  '

Set-TextValue $ws.Range("E13") '  -0.27%  '

Set-TextValue $ws.Range("D14") '2.742.67'
Set-TextValue $ws.Range("E14") '  +3.33%  '

Set-TextValue $ws.Range("D15") '2.359.79'
Set-TextValue $ws.Range("E15") '  +2.58%  '

Set-TextValue $ws.Range("D16") '0.826'
Set-TextValue $ws.Range("E16") '  +0.27%  '

Set-TextValue $ws.Range("D17") '13.77'
Set-TextValue $ws.Range("E17") '  -2.12%  '

Set-TextValue $ws.Range("D18") '45.972.32'
Set-TextValue $ws.Range("E18") '  -1.74%  '

Set-TextValue $ws.Range("D19") '12.66'
Set-TextValue $ws.Range("E19") '  -5.51%  '

Set-TextValue $ws.Range("D20") '0.0₃0951'
Set-TextValue $ws.Range("E20") '  -0.19%  '

Set-TextValue $ws.Range("D21") '6.05'
Set-TextValue $ws.Range("E21") '  -1.30%  '

Set-TextValue $ws.Range("D22") '66.79'
Set-TextValue $ws.Range("E22") '  -0.13%  '

Set-TextValue $ws.Range("D23") '243.78'
Set-TextValue $ws.Range("E23") '  -2.02%  '

Set-TextValue $ws.Range("D24") '2.79'
Set-TextValue $ws.Range("E24") '  -5.62%  '

Set-TextValue $ws.Range("E25") '  +0.11%  '

Set-TextValue $ws.Range("D26") '1.93'
Set-TextValue $ws.Range("E26") '  -2.21%  '

Set-TextValue $ws.Range("D27") '39.96'
Set-TextValue $ws.Range("E27") '  -10.37%  '

Set-TextValue $ws.Range("D28") '2.21'
Set-TextValue $ws.Range("E28") '  -3.06%  '

Set-TextValue $ws.Range("D29") '9.74'
Set-TextValue $ws.Range("E29") '  -2.21%  '

Set-TextValue $ws.Range("D30") '20.94'
Set-TextValue $ws.Range("E30") '  +3.87%  '

Set-TextValue $ws.Range("E31") '  +17.56%  '

Set-TextValue $ws.Range("E32") '  +7.02%  '

Set-TextValue $ws.Range("E33") '  -4.42%  '

Set-TextValue $ws.Range("D34") '147.01'
Set-TextValue $ws.Range("E34") '  -0.45%  '

Set-TextValue $ws.Range("D35") '0.0773'
Set-TextValue $ws.Range("E35") '  -3.63%  '

Set-TextValue $ws.Range("E36") '  -0.59%  '

Set-TextValue $ws.Range("D37") '1.92'
Set-TextValue $ws.Range("E37") '  +5.86%  '

Set-TextValue $ws.Range("E38") '  -2.83%  '

Set-TextValue $ws.Range("D39") '15.08'
Set-TextValue $ws.Range("E39") '  -6.60%  '

Set-TextValue $ws.Range("D40") '3.88'
Set-TextValue $ws.Range("E40") '  -4.08%  '

Set-TextValue $ws.Range("E41") '  -2.17%  '

Set-TextValue $ws.Range("E42") '  -8.24%  '

Set-TextValue $ws.Range("D43") '1.934.52'
Set-TextValue $ws.Range("E43") '  +3.46%  '

Set-TextValue $ws.Range("E44") '  +0.07%  '

Set-TextValue $ws.Range("D45") '92.15'
Set-TextValue $ws.Range("E45") '  +4.43%  '

Set-TextValue $ws.Range("E46") '  -9.74%  '

Set-TextValue $ws.Range("E47") '  +5.50%  '

Set-TextValue $ws.Range("E48") '  -5.47%  '

Set-TextValue $ws.Range("D49") '98.97'
Set-TextValue $ws.Range("E49") '  +1.65%  '

Set-TextValue $ws.Range("D50") '2.613.23'
Set-TextValue $ws.Range("E50") '  +3.32%  '

Set-TextValue $ws.Range("D51") '68.90'
Set-TextValue $ws.Range("E51") '  -7.67%  '
